# Actualización automática 2025-06-11 11:00:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("N31").Value = 1520.03
$wsGrupo.Range("D34").Value = 814.08
$wsGrupo.Range("D58").Value = "1 de 56"
$wsGrupo.Range("N58").Value = "1 de 56"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F31").Value = 1520.03
$wsMensual.Range("F34").Value = 814.08
$wsMensual.Range("F58").Value = 4776.17

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D3").Value = 814.08
$wsCumpl.Range("E3").Value = 4178.1032
$wsCumpl.Range("F3").Value = 0.1630709385825424

$wsCumpl.Range("D18").Value = 1520.03
$wsCumpl.Range("E18").Value = 2609.97
$wsCumpl.Range("F18").Value = 0.368046004842615

$wsCumpl.Range("D19").Value = 4776.17
$wsCumpl.Range("E19").Value = 50633.53560036206
$wsCumpl.Range("F19").Value = 0.08619735384352575
